# Insert a new "most recent period" column before column D on the MIK sheet,
# shifting the existing D:K data right to E:L, then fill the new column D
# with the latest financial figures (matching the author's update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MIK")

# Insert a new column at D; existing D:K shift to E:L (formats shift with
# the cells, bringing each column's own number formatting along).
$ws.Columns("D:D").Insert()

# The freshly inserted column picks up the formatting of the column to its
# left (C) by default. Re-apply the per-row number formats from the shifted
# former-D column, now sitting in E, so the new D matches the rest of its row
# (date format for the header rows, thousands format everywhere else).
# Copy block-by-block (matching the sheet's existing used-row blocks) so we
# don't materialize cells on blank separator rows (36, 78) that have no data.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the new period's values.
$ws.Range("D7").Value = 43498
$ws.Range("D8").Value = 5271900
$ws.Range("D9").Value = 3248300
$ws.Range("D10").Value = 2023700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 106100
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 4710200
$ws.Range("D18").Value = 561800
$ws.Range("D20").Value = 2400
$ws.Range("D21").Value = 688400
$ws.Range("D22").Value = 147100
$ws.Range("D23").Value = 417100
$ws.Range("D24").Value = 96500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 320500
$ws.Range("D27").Value = 320000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -1000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2400
$ws.Range("D33").Value = 319000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 319000
$ws.Range("D38").Value = 43498
$ws.Range("D41").Value = 245900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 62300
$ws.Range("D44").Value = 1108700
$ws.Range("D45").Value = 98700
$ws.Range("D46").Value = 1515500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 439100
$ws.Range("D49").Value = 129300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 44400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2128300
$ws.Range("D57").Value = 485000
$ws.Range("D58").Value = 24900
$ws.Range("D59").Value = 422600
$ws.Range("D60").Value = 932600
$ws.Range("D61").Value = 2681000
$ws.Range("D62").Value = 141000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 3754500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -1628200
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = -1626200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43498
$ws.Range("D81").Value = 319000
$ws.Range("D83").Value = 124300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 444300
$ws.Range("D91").Value = -145400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -145400
$ws.Range("D96").Value = -300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -478900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -180000

$wb.Save()
